# Applies the "expermits todos no convexos menos el 5to" edit:
# updates the recorded non-convex-experiment coefficients/results across
# several sheets of the workbook.
#
# Most of the target cells hold text that LOOKS like a plain number
# (e.g. "0.79", "-0.0", "10.0"). Those sheets store every value as a
# shared string (t="s"), so a naive `.Value = "0.79"` would get silently
# re-typed as a numeric cell by Excel's input parser. To keep the cell a
# text value (matching the original authoring) we briefly force the
# cell to Text number-format, assign the literal string, then restore
# the cell's format back to General/Normal so no visible formatting
# change is left behind.

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Restricciones_del_follower
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2 (J_0_L0_v)
$ws.Range("A2").Value = "7.994501225455026 - 2x_1 + 1.12805778384666y_1 - 0.17512553170551792y_2"
Set-TextValue $ws.Range("B2") "-5.494501225455026"
Set-TextValue $ws.Range("D2") "0.79"
Set-TextValue $ws.Range("E2") "10.0"
Set-TextValue $ws.Range("F2") "0"

# Row 3 (J_0_LP_v)
$ws.Range("A3").Value = "3.4975717805348214 + x_1 - 3x_2 - 0.2730288181213387y_1 - 0.758694351969956y_2"
Set-TextValue $ws.Range("B3") "-5.497571780534821"
Set-TextValue $ws.Range("D3") "0.09"
Set-TextValue $ws.Range("E3") "-3.2"
Set-TextValue $ws.Range("F3") "-0.0"

# Row 4 (J_Ne_L0_v)
$ws.Range("A4").Value = "-12.317746389472813 + x_1 + x_2 + 0.13628043375465734y_1 + 0.8778400418589355y_2"
Set-TextValue $ws.Range("B4") "10.267746389472812"
Set-TextValue $ws.Range("D4") "0.54"
Set-TextValue $ws.Range("E4") "0.3"
Set-TextValue $ws.Range("F4") "1.6"

# ---------------------------------------------------------------
# Punto_modificado
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "6.65"
Set-TextValue $ws.Range("B2") "2.1"
Set-TextValue $ws.Range("C2") "5.2"
Set-TextValue $ws.Range("D2") "3.2"

# ---------------------------------------------------------------
# Vector_bf  (sheet index 5 -- name lookup is case-insensitive, and
# this workbook also has a "Vector_BF" sheet, so address by index to
# avoid ambiguity between the two names)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "3.059815510164544"
Set-TextValue $ws.Range("A3") "-1.26740196087917"

# ---------------------------------------------------------------
# Vector_BF (sheet index 6)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "24.9"
Set-TextValue $ws.Range("A3") "-10.900000000000002"
Set-TextValue $ws.Range("A4") "-12.695154186581279"
Set-TextValue $ws.Range("A5") "-0.9399186218063611"

# ---------------------------------------------------------------
# Vector_Alpha (these two cells are genuine numbers, not text)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 0.27629937218188366
$ws.Range("A3").Value = 1.7797613770322578
